$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.096.39"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.788.95"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.36"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.14"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.47"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0670"
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0923"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.043.86"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.788.65"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.626"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.045.15"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.19"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.93"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "252.76"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0745"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.40"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.21"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.54"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.57"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.01"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.86"
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.505.08"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.639"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.60"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0514"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.08"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.940.24"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.75"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.87"
$ws.Range("E50").Value = "  +14.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.57"
$ws.Range("E51").Value = "  -6.10%  "
